# Add 2022-Q4 data:
#  - insert a new worksheet "2022-Q4" right after "总计", pushing
#    "2022-Q3" and "2021-Q4" one position to the right
#  - populate the new sheet with the 2022-Q4 fund holdings table
#  - update the "总计" (summary) sheet with a new row for 2022-Q4 and
#    shift the existing 2022-Q3 / 2021-Q4 rows down by one

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "2022-Q4" worksheet right after "总计" ---------------
$totalSheetForAdd = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($null, $totalSheetForAdd)
$newSheet.Name = "2022-Q4"

# Re-fetch sheet references now that the sheet collection has changed, so we
# are not holding on to any reference captured before the insertion.
$totalSheet = $wb.Worksheets.Item("总计")
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$newSheet = $wb.Worksheets.Item("2022-Q4")

# Bring over the header row (and its styling) plus the A-column index style
# used by the other quarterly sheets, so the new sheet matches their look.
$q3Sheet.Range("A1:H1").Copy($newSheet.Range("A1"))
$q3Sheet.Range("A2:A3").Copy($newSheet.Range("A2"))

# --- 2. Fill in the 2022-Q4 fund data ---------------------------------------
# Columns B-G hold numeric-looking text (fund codes, percentages, etc.) in
# the source data, so force a text number format before assigning them to
# avoid Excel re-interpreting them as numbers (which would e.g. drop the
# leading zero of a fund code or the trailing zero of a percentage).
$newSheet.Range("B2:G3").NumberFormat = "@"

$newSheet.Range("B2").Value = "011807"
$newSheet.Range("C2").Value = "平安研究精选混合A"
$newSheet.Range("D2").Value = "1.23"
$newSheet.Range("E2").Value = "86.30"
$newSheet.Range("F2").Value = "2.16"
$newSheet.Range("G2").Value = "0.0266"
$newSheet.Range("H2").Value = 10

$newSheet.Range("B3").Value = "011808"
$newSheet.Range("C3").Value = "平安研究精选混合C"
$newSheet.Range("D3").Value = "1.05"
$newSheet.Range("E3").Value = "86.30"
$newSheet.Range("F3").Value = "2.16"
$newSheet.Range("G3").Value = "0.0227"
$newSheet.Range("H3").Value = 10

# --- 3. Update the "总计" summary sheet -------------------------------------
# Give the new row-4 index cell (A4) the same style as the existing A3 index
# cell before shifting values, so it matches the formatting used by A2/A3.
$totalSheet.Range("A3").Copy($totalSheet.Range("A4"))

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.05

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q3"
$totalSheet.Range("C3").Value = 4
$totalSheet.Range("D3").Value = 0.4

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q4"
$totalSheet.Range("C4").Value = 2
$totalSheet.Range("D4").Value = 0.16
